$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.032.47'
$ws.Range("E2").Value = '  -4.39%  '
$ws.Range("D3").Value = '2.971.48'
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.66%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.516'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.26%  '
$ws.Range("D9").Value = '2.963.78'
$ws.Range("E9").Value = '  -1.59%  '
$ws.Range("E10").Value = '  -2.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.90%  '
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000226'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '3.461.37'
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.13%  '
$ws.Range("D18").Value = '2.967.76'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").Value = '57.965.31'
$ws.Range("E19").Value = '  -4.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '422.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.690'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0996'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("E35").Value = '  -1.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.945'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("D37").Value = '0.0₃0701'
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.84%  '
$ws.Range("E40").Value = '  +4.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0353'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.07%  '
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '379.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("D44").Value = '2.710.83'
$ws.Range("E44").Value = '  +1.71%  '
$ws.Range("E46").Value = '  +2.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.111'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.54%  '
